$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells stay text (avoid Excel auto-converting numeric-looking strings)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.759.86"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.282.49"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.71"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.636"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.59"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +6.54%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.638"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.70"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.22%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.45"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.622.23"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.02"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.03%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.280.86"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.641.02"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.24"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.45"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "237.14"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.16"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.86"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.31"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.41"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.77"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.04"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0888"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +12.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.33"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.126"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.40"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.05%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.09%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.15%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.91%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.69"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +9.69%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.90"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.95%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.04"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "61.24"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "105.33"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +11.66%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.23"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.84%  "
